# Atualizado por script em 03-11-2023 14:45
#
# This script:
#  1) Swaps the match data (columns F:V, excluding the shared timestamp
#     columns K/O/S) between rows 57 and 58, which were reordered.
#  2) Appends a new match row (row 74) for Sesvete vs Croatia Zmijavci,
#     copying the formatting from the last existing data row (73) and
#     then filling in its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap rows 57 and 58 (data columns only; A-E and K/O/S untouched)
# ---------------------------------------------------------------

# Row 57 becomes what used to be row 58 (Solin vs Jarun)
$ws.Cells.Item(57,6).Value  = "Solin"
$ws.Cells.Item(57,7).Value  = 1
$ws.Cells.Item(57,8).Value  = "Jarun"
$ws.Cells.Item(57,9).Value  = 1
$ws.Cells.Item(57,10).Value = 1.85
$ws.Cells.Item(57,12).Value = 1.88
$ws.Cells.Item(57,13).Value = "14/10/2023 14:51"
$ws.Cells.Item(57,14).Value = 3.54
$ws.Cells.Item(57,16).Value = 3.62
$ws.Cells.Item(57,17).Value = "14/10/2023 14:51"
$ws.Cells.Item(57,18).Value = 3.41
$ws.Cells.Item(57,20).Value = 3.79
$ws.Cells.Item(57,21).Value = "14/10/2023 14:51"
$ws.Cells.Item(57,22).Value = "https://www.betexplorer.com/football/croatia/prva-nl/solin-jarun/6BAb7QlU/"

# Row 58 becomes what used to be row 57 (Zrinski Jurjevac vs Croatia Zmijavci)
$ws.Cells.Item(58,6).Value  = "Zrinski Jurjevac"
$ws.Cells.Item(58,7).Value  = 4
$ws.Cells.Item(58,8).Value  = "Croatia Zmijavci"
$ws.Cells.Item(58,9).Value  = 0
$ws.Cells.Item(58,10).Value = 1.68
$ws.Cells.Item(58,12).Value = 1.58
$ws.Cells.Item(58,13).Value = "14/10/2023 14:53"
$ws.Cells.Item(58,14).Value = 3.61
$ws.Cells.Item(58,16).Value = 3.86
$ws.Cells.Item(58,17).Value = "14/10/2023 14:53"
$ws.Cells.Item(58,18).Value = 4.26
$ws.Cells.Item(58,20).Value = 5.57
$ws.Cells.Item(58,21).Value = "14/10/2023 14:53"
$ws.Cells.Item(58,22).Value = "https://www.betexplorer.com/football/croatia/prva-nl/zrinski-jurjevac-croatia-zmijavci/dzj8RTZo/"

# ---------------------------------------------------------------
# 2) Append new row 74 (Sesvete vs Croatia Zmijavci)
# ---------------------------------------------------------------

# Clone formatting from the previous last row (73) onto the new row (74)
$ws.Range("A73:V73").Copy()
$ws.Range("A74:V74").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(74,1).Value  = 73
$ws.Cells.Item(74,2).Value  = "croatia"
$ws.Cells.Item(74,3).Value  = "prva-nl"
$ws.Cells.Item(74,4).Value  = "2023-2024"
$ws.Cells.Item(74,5).Value  = 45233.58333333334
$ws.Cells.Item(74,6).Value  = "Sesvete"
$ws.Cells.Item(74,7).Value  = 1
$ws.Cells.Item(74,8).Value  = "Croatia Zmijavci"
$ws.Cells.Item(74,9).Value  = 0
$ws.Cells.Item(74,10).Value = 1.95
$ws.Cells.Item(74,11).Value = "02/11/2023 02:12"
$ws.Cells.Item(74,12).Value = 1.7
$ws.Cells.Item(74,13).Value = "03/11/2023 13:51"
$ws.Cells.Item(74,14).Value = 3.32
$ws.Cells.Item(74,15).Value = "02/11/2023 02:12"
$ws.Cells.Item(74,16).Value = 3.2
$ws.Cells.Item(74,17).Value = "03/11/2023 13:51"
$ws.Cells.Item(74,18).Value = 3.31
$ws.Cells.Item(74,19).Value = "02/11/2023 02:12"
$ws.Cells.Item(74,20).Value = 5.85
$ws.Cells.Item(74,21).Value = "03/11/2023 13:51"
$ws.Cells.Item(74,22).Value = "https://www.betexplorer.com/football/croatia/prva-nl/sesvete-croatia-zmijavci/U7pXHG2k/"
